# Add a new "Compact List" paragraph style, cloned from the existing
# "Compact" style (same base style, quick-style flag, and before/after
# paragraph spacing of 36 twips == 1.8 pt).
$d = $word.ActiveDocument

$newStyle = $d.Styles.Add("Compact List", 1)   # 1 = wdStyleTypeParagraph
$newStyle.BaseStyle = "BodyText"
$newStyle.QuickStyle = $true
$newStyle.ParagraphFormat.SpaceBefore = 1.8
$newStyle.ParagraphFormat.SpaceAfter = 1.8
